# #327 Ajout des profils d'acces
#
# 1. Bump the "Date" metadata value on the Metadata sheet.
# 2. On the Elements sheet, the two "Mapping" columns (AK = RIM Mapping,
#    AL = Spécification métier vers l'extension ROR ContactTelecomUsaget)
#    swap places: the "Spécification métier..." mapping now comes first
#    (column AK) and "RIM Mapping" moves to column AL. Swap the cell
#    contents of the two columns row-by-row so the data follows the header.

$wb = $excel.ActiveWorkbook

# --- 1. Update the "Date" metadata value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2. Swap the content of columns AK (37) and AL (38) on "Elements" ---
$els = $wb.Worksheets.Item("Elements")

$lastRow = $els.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $els.Cells.Item($r, 37)
    $alCell = $els.Cells.Item($r, 38)
    $akVal = $akCell.Value()
    $alVal = $alCell.Value()
    if ($akVal -ne $alVal) {
        $akCell.Value = $alVal
        $alCell.Value = $akVal
    }
}
